$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction in TSR formula: the TSR column (F) previously computed
# SUM(1, -D, -E) for rows 4, 6 and 8. Update it to the weighted formula
# (5 * SUM(1, -E) + 4 * SUM(1, -D)) / 9
$ws.Range("F4").Formula = "= (5 * SUM(1, -E4) + 4 * SUM(1, -D4)) / 9"
$ws.Range("F6").Formula = "= (5 * SUM(1, -E6) + 4 * SUM(1, -D6)) / 9"
$ws.Range("F8").Formula = "= (5 * SUM(1, -E8) + 4 * SUM(1, -D8)) / 9"

# Update the active selection left in the sheet
[void]$ws.Range("F13").Select()
